$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'52.416.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.74%  "

# Row 3
$ws.Range("D3").Value = "'2.923.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.85%  "

# Row 4
$ws.Range("E4").Value = "  +0.03%  "

# Row 5
$ws.Range("D5").Value = "'352.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "

# Row 6
$ws.Range("D6").Value = "'112.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.64%  "

# Row 7
$ws.Range("E7").Value = "  +1.42%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.628"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.33%  "

# Row 10
$ws.Range("D10").Value = "'40.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.30%  "

# Row 11
$ws.Range("D11").Value = "'0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.93%  "

# Row 12
$ws.Range("E12").Value = "  +0.50%  "

# Row 13
$ws.Range("D13").Value = "'20.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.32%  "

# Row 14
$ws.Range("D14").Value = "'7.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.82%  "

# Row 15
$ws.Range("D15").Value = "'3.381.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.90%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'2.910.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.30%  "

# Row 17
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.993"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.71%  "

# Row 18
$ws.Range("D18").Value = "'52.438.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.93%  "

# Row 19
$ws.Range("D19").Value = "'14.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.42%  "

# Row 20
$ws.Range("D20").Value = "'7.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.57%  "

# Row 21
$ws.Range("E21").Value = "  +5.76%  "

# Row 22
$ws.Range("D22").Value = "'0.0₃0982"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.23%  "

# Row 23
$ws.Range("D23").Value = "'71.18"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.31%  "

# Row 24
$ws.Range("D24").Value = "'271.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.66%  "

# Row 25
$ws.Range("D25").Value = "'2.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.52%  "

# Row 26
$ws.Range("D26").Value = "'26.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.81%  "

# Row 27
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "

# Row 28
$ws.Range("E28").Value = "  +0.11%  "

# Row 29
$ws.Range("D29").Value = "'10.65"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.38%  "

# Row 30
$ws.Range("D30").Value = "'38.11"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.14%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'6.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.69%  "

# Row 32
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +13.86%  "

# Row 33
$ws.Range("D33").Value = "'6.17"
$ws.Range("D33").Style = "Normal"

# Row 34
$ws.Range("D34").Value = "'53.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.58%  "

# Row 35
$ws.Range("D35").Value = "'0.0937"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.02%  "

# Row 36
$ws.Range("D36").Value = "'0.0451"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.95%  "

# Row 37
$ws.Range("E37").Value = "  -0.15%  "

# Row 38
$ws.Range("D38").Value = "'3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.57%  "

# Row 39
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "'2.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.03%  "

# Row 40
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'18.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.18%  "

# Row 41
$ws.Range("D41").Value = "'2.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.43%  "

# Row 42
$ws.Range("D42").Value = "'24.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.57%  "

# Row 43
$ws.Range("E43").Value = "  +2.13%  "

# Row 44
$ws.Range("D44").Value = "'122.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.53%  "

# Row 45
$ws.Range("E45").Value = "  +1.02%  "

# Row 46
$ws.Range("D46").Value = "'3.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.92%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'2.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.03%  "

# Row 48
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "'2.221.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.78%  "

# Row 49
$ws.Range("D49").Value = "'0.265"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +25.76%  "

# Row 50
$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D50").Value = "'0.0337"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.47%  "

# Row 51
$ws.Range("B51").Value = "SEI"
$ws.Range("C51").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D51").Value = "'0.966"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.92%  "
